$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-5 (Q0-Q3) with new values
$ws.Range("B2").Value = 0.6158830556516814
$ws.Range("C2").Value = 2.146649694053442
$ws.Range("D2").Value = 6.552896975861514
$ws.Range("E2").Value = 2.559862686915357
$ws.Range("F2").Value = 2.578463921477837
$ws.Range("G2").Value = 14

$ws.Range("B3").Value = 0.6396189027044379
$ws.Range("C3").Value = 1.873773123107626
$ws.Range("D3").Value = 5.190755540198163
$ws.Range("E3").Value = 2.278322966613417
$ws.Range("F3").Value = 2.30498088676128
$ws.Range("G3").Value = 10

$ws.Range("B4").Value = 1.171938100156769
$ws.Range("C4").Value = 1.236930977686064
$ws.Range("D4").Value = 2.577565538177816
$ws.Range("E4").Value = 1.605479846705594
$ws.Range("F4").Value = 1.20206154297295
$ws.Range("G4").Value = 6

$ws.Range("B5").Value = 0.3821426306726097
$ws.Range("C5").Value = 0.3821426306726097
$ws.Range("D5").Value = 0.206649130868104
$ws.Range("E5").Value = 0.4545867693500373
$ws.Range("F5").Value = 0.3481842635465348
$ws.Range("G5").Value = 2

# Remove rows 6-9 (Q4-Q7), which also drops their shared strings
$ws.Range("A6:G9").Delete()
